$wb = $excel.ActiveWorkbook

# --- Add the new worksheet as the last tab (right after "Kevin's Tests") ---
$kevin = $wb.Worksheets.Item("Kevin's Tests")
$new = $wb.Worksheets.Add($null, $kevin)
$new.Name = "Marlene's Tests"

# --- Header row (row 1) ---
$new.Range("A1").Value = "Method"
$new.Range("B1").Value = " 1 -> 2"
$new.Range("C1").Value = " 1 -> 3"
$new.Range("D1").Value = " 1 -> 4"
$new.Range("E1").Value = " 2 -> 3"
$new.Range("F1").Value = " 2 -> 4"
$new.Range("G1").Value = " 3 -> 4"
$new.Range("H1").Value = " Time"

# Re-use the existing bold+centered header style (style index 1) from the
# "Kevin's Tests" sheet instead of constructing new font/alignment formatting
# (keeps the shared style table untouched).
$kevin.Range("A1:H1").Copy() | Out-Null
$new.Range("A1:H1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Row 2: Levenshtein's Distance ---
$new.Range("A2").Value = "Levenshtein's Distance"
$new.Range("B2").Value = 28
$new.Range("C2").Value = 19
$new.Range("D2").Value = 32
$new.Range("E2").Value = 11
$new.Range("F2").Value = 20
$new.Range("G2").Value = 28
$new.Range("H2").Value = 0.0311574

# --- Row 3: Jaro-Winkler Distance ---
$new.Range("A3").Value = "Jaro-Winkler Distance"
$new.Range("B3").Value = 0.646719
$new.Range("C3").Value = 0.706742
$new.Range("D3").Value = 0.476475
$new.Range("E3").Value = 0.854641
$new.Range("F3").Value = 0.849474
$new.Range("G3").Value = 0.60762
$new.Range("H3").Value = 0.0155787

# --- Row 4: Hunt-McIlroy Distance ---
$new.Range("A4").Value = "Hunt-McIlroy Distance"
$new.Range("B4").Value = 37
$new.Range("C4").Value = 37
$new.Range("D4").Value = 18
$new.Range("E4").Value = 44
$new.Range("F4").Value = 33
$new.Range("G4").Value = 30
$new.Range("H4").Value = 0.0311574

# --- Row 5: Needleman-Wunsch Distance ---
$new.Range("A5").Value = "Needleman-Wunsch Distance"
$new.Range("B5").Value = -21
$new.Range("C5").Value = -6
$new.Range("D5").Value = -47
$new.Range("E5").Value = 18
$new.Range("F5").Value = -9
$new.Range("G5").Value = -12
$new.Range("H5").Value = 0.0467361

# Column A best-fits the longest label ("Needleman-Wunsch Distance"), same as
# the analogous column on "Greg's tests".
$new.Columns.Item(1).AutoFit() | Out-Null

# Selection / active-cell bookkeeping matching the authored sheet.
$new.Range("L9").Select()

# Make the newly authored sheet the active tab (mirrors tabSelected moving to
# this sheet and off of "Result2").
$new.Activate()
